# Simulated Wild Card round and logged it
$wb = $excel.ActiveWorkbook

# --- OFF sheet: update row 2 (H) stats with Wild Card round numbers ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 339
$wsOff.Range("C2").Value = 231
$wsOff.Range("D2").Value = 165
$wsOff.Range("E2").Value = 64

# --- DEF sheet: update row 2 (H) stats with Wild Card round numbers ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 388
$wsDef.Range("C2").Value = 283
$wsDef.Range("D2").Value = 103
$wsDef.Range("E2").Value = 51
$wsDef.Range("G2").Value = 6
